$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "MIE(8.481922935287269, 4.237989323052444, -9.965146422450527, 11.139222298904546)"
$ws.Range("C2").Value = "NIG(1.0513918225868284, 0.7760280958142305, 4.48075767548238, 6.098304193269812)"
$ws.Range("D2").Value = "NCT(2.3734245980568076, 1.3136492526689105, -0.16244699254479597, 2.329310234945339)"
$ws.Range("E2").Value = "NCT(3.332378016842687, 1.6321949196154626, 1.0256679312456434, 4.620897009046477)"
